$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force-text cells: numeric-looking Price strings that must stay text (matches source data formatting)
$textCells = @("D5", "D10", "D11", "D16", "D17", "D19", "D25", "D27", "D37", "D38", "D39", "D41", "D42", "D44", "D47")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply updated cell values
$ws.Range("D2").Value = "26.976.31"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").Value = "1.677.60"
$ws.Range("E3").Value = "  +0.40%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "214.95"
$ws.Range("E5").Value = "  -0.51%  "
$ws.Range("E6").Value = "  +1.03%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("E8").Value = "  -0.29%  "
$ws.Range("E9").Value = "  +0.24%  "
$ws.Range("D10").Value = "20.31"
$ws.Range("E10").Value = "  +1.07%  "
$ws.Range("D11").Value = "0.0887"
$ws.Range("E11").Value = "  -0.79%  "
$ws.Range("D12").Value = "1.913.82"
$ws.Range("E12").Value = "  +0.40%  "
$ws.Range("D13").Value = "1.692.10"
$ws.Range("E13").Value = "  +1.36%  "
$ws.Range("E14").Value = "  +0.14%  "
$ws.Range("E15").Value = "  +1.51%  "
$ws.Range("D16").Value = "65.69"
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").Value = "8.20"
$ws.Range("E17").Value = "  +6.06%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "27.001.94"
$ws.Range("E18").Value = "  -0.12%  "
$ws.Range("D19").Value = "235.47"
$ws.Range("E19").Value = "  +0.07%  "
$ws.Range("D20").Value = "0.0₃0733"
$ws.Range("E20").Value = "  -0.35%  "
$ws.Range("E21").Value = "  +0.11%  "
$ws.Range("E22").Value = "  -0.21%  "
$ws.Range("E23").Value = "  -0.61%  "
$ws.Range("E24").Value = "  -2.67%  "
$ws.Range("D25").Value = "146.12"
$ws.Range("E25").Value = "  +0.42%  "
$ws.Range("E26").Value = "  +0.68%  "
$ws.Range("D27").Value = "16.06"
$ws.Range("E27").Value = "  +1.16%  "
$ws.Range("E28").Value = "  -1.49%  "
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("E30").Value = "  -0.10%  "
$ws.Range("E31").Value = "  -0.75%  "
$ws.Range("E32").Value = "  -0.20%  "
$ws.Range("D33").Value = "1.478.26"
$ws.Range("E33").Value = "  +1.88%  "
$ws.Range("E34").Value = "  +1.09%  "
$ws.Range("E35").Value = "  +4.86%  "
$ws.Range("E36").Value = "  +0.17%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "0.0175"
$ws.Range("E37").Value = "  +2.91%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "0.582"
$ws.Range("E38").Value = "  +2.03%  "
$ws.Range("D39").Value = "0.901"
$ws.Range("E39").Value = "  +0.90%  "
$ws.Range("E40").Value = "  -3.95%  "
$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.15%  "
$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").Value = "1.02"
$ws.Range("E42").Value = "  +0.46%  "
$ws.Range("E43").Value = "  +0.87%  "
$ws.Range("D44").Value = "67.40"
$ws.Range("E44").Value = "  +2.49%  "
$ws.Range("D45").Value = "1.818.23"
$ws.Range("E45").Value = "  +0.14%  "
$ws.Range("D47").Value = "90.42"
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("E48").Value = "  +0.84%  "
$ws.Range("E49").Value = "  -0.62%  "
$ws.Range("E50").Value = "  +1.35%  "
$ws.Range("E51").Value = "  +0.10%  "

# Restore default style for cells we force-formatted as text (keeps formatting identical to source)
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).Style = "Normal"
}
